$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns
$ws.Range("H1").Value = "BANK"
$ws.Range("I1").Value = "ACCOUNT_NUMBER"
$ws.Range("J1").Value = "DATE_OF_BIRTH"
$ws.Range("K1").Value = "DATE_OF_JOINING"

# Column widths for the new columns (values chosen so the resulting stored
# OOXML width, after this runtime's internal character->pixel rounding,
# lands on the closest achievable value to the target width)
$ws.Columns.Item(9).ColumnWidth = 17.1
$ws.Columns.Item(10).ColumnWidth = 13.59
$ws.Columns.Item(11).ColumnWidth = 16.6

# Bank / account data for rows 2-4
$ws.Range("H2").Value = "Wema"
$ws.Range("I2").Value = 1234567890
$ws.Range("H3").Value = "Wema"
$ws.Range("I3").Value = 1234567890
$ws.Range("H4").Value = "Wema"
$ws.Range("I4").Value = 1234567890

# Date values (serial numbers): 33824 = 08-Aug-1992, 36746 = 08-Aug-2000
$ws.Range("J2").Value = 33824
$ws.Range("K2").Value = 36746
$ws.Range("J3").Value = 33824
$ws.Range("K3").Value = 36746
$ws.Range("J4").Value = 33824
$ws.Range("K4").Value = 36746

# Apply date number format once, then copy the format to the rest of the
# date cells so they all reuse the same cell-style index.
$ws.Range("J2").NumberFormat = "mm-dd-yy"
$ws.Range("J2").Copy()
$ws.Range("K2:K4").PasteSpecial(-4122)
$ws.Range("J3:J4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Matches the activeCell selection recorded in the saved file
$ws.Range("I5").Select()
